# Apply cryptos.xlsx symbol-list update (Tue Jan  3 19:20:39 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) hold numeric-looking text values
# (stored as text in the sheet, e.g. "245.00", "-0.60%"). Force the
# cells to Text format before assigning so Excel keeps them as literal
# strings (with exact formatting / trailing zeros) instead of coercing
# them into numbers or percentage fractions.
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.50%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.28%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.234"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.14%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05698"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.55%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.611"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.26%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8510"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.63%"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8596"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.84%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1366"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.33%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07026"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.56%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03136"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.79%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09231"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.64%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001534"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.36%"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005955"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.53%"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.488"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.00%"

# Row 17
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.174"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.33%"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005988"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-94.17%"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.42%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03266"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-5.34%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1298"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.45%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.485"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.02%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04086"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.95%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1379"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.08%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001224"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.63%"

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-17.57%"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.84%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03764"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.32%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1062"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.03%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003735"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.12%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002300"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.46%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009149"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.32%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005276"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.07%"

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.01%"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1151"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "78.00%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002439"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.23%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.01%"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
